$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10206857
$ws.Range("I19").Value = 6708721
$ws.Range("K19").Value = 6708721
$ws.Range("M19").Value = -6708546

$ws.Range("H33").Value = 90.5
$ws.Range("I33").Value = 98.44444
$ws.Range("K33").Value = 98.44444
$ws.Range("M33").Value = 130.55556

$ws.Range("H82").Value = 4800147
$ws.Range("I82").Value = 4800147
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 14400441
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -14400035
$ws.Range("N82").Value = ""

$ws.Range("H85").Value = 4800147
$ws.Range("I85").Value = 4800147
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 14400441
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -14399037
$ws.Range("N85").Value = ""

$ws.Range("H116").Value = 1510248
$ws.Range("I116").Value = 1604457.2
$ws.Range("J116").Value = 2900
$ws.Range("K116").Value = 1604457.2
$ws.Range("L116").Value = 2900
$ws.Range("M116").Value = -1601015.2
$ws.Range("N116").Value = -9784

$ws.Range("H132").Value = 1969.7455
$ws.Range("I132").Value = 1985.3077
$ws.Range("K132").Value = 5955.9231
$ws.Range("M132").Value = -3425.9231

$ws.Range("H137").Value = 8369114.5
$ws.Range("I137").Value = 13890608
$ws.Range("J137").Value = 86874.5
$ws.Range("K137").Value = 41671824
$ws.Range("L137").Value = 260623.5
$ws.Range("M137").Value = -41669274
$ws.Range("N137").Value = -265723.5

$ws.Range("H138").Value = 2621.9565
$ws.Range("I138").Value = 1465.9
$ws.Range("J138").Value = 3511.2307
$ws.Range("K138").Value = 4397.700000000001
$ws.Range("L138").Value = 10533.6921
$ws.Range("M138").Value = 742.2999999999993
$ws.Range("N138").Value = -20813.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 544.0526
$ws.Range("I2").Value = 472
$ws.Range("J2").Value = 1156.5
$ws.Range("K2").Value = 472
$ws.Range("L2").Value = 1156.5
$ws.Range("M2").Value = -359
$ws.Range("N2").Value = -1382.5

$ws.Range("H32").Value = 18187120
$ws.Range("I32").Value = 20002112
$ws.Range("J32").Value = 37199.8
$ws.Range("K32").Value = 20002112
$ws.Range("L32").Value = 37199.8
$ws.Range("M32").Value = -20001825
$ws.Range("N32").Value = -37773.8

$ws.Range("H61").Value = 1138.8406
$ws.Range("I61").Value = 1226.7675
$ws.Range("J61").Value = 993.4231
$ws.Range("K61").Value = 1226.7675
$ws.Range("L61").Value = 993.4231
$ws.Range("M61").Value = -1014.7675
$ws.Range("N61").Value = -1417.4231

$ws.Range("H116").Value = 544.0526
$ws.Range("I116").Value = 472
$ws.Range("J116").Value = 1156.5
$ws.Range("K116").Value = 472
$ws.Range("L116").Value = 1156.5
$ws.Range("M116").Value = 1822
$ws.Range("N116").Value = -5744.5

$ws.Range("H136").Value = 1138.8406
$ws.Range("I136").Value = 1226.7675
$ws.Range("J136").Value = 993.4231
$ws.Range("K136").Value = 3680.3025
$ws.Range("L136").Value = 2980.2693
$ws.Range("M136").Value = -1130.3025
$ws.Range("N136").Value = -8080.2693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 544.0526
$ws.Range("I3").Value = 472
$ws.Range("J3").Value = 1156.5
$ws.Range("K3").Value = 472
$ws.Range("L3").Value = 1156.5
$ws.Range("M3").Value = -358
$ws.Range("N3").Value = -1384.5

$ws.Range("H22").Value = 102
$ws.Range("J22").Value = 102
$ws.Range("L22").Value = 102
$ws.Range("N22").Value = -448

$ws.Range("H25").Value = 1225.5555
$ws.Range("I25").Value = 1225.5555
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1225.5555
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -990.5554999999999
$ws.Range("N25").Value = ""

$ws.Range("H99").Value = 1721.2051
$ws.Range("I99").Value = 1139.2858
$ws.Range("J99").Value = 3202.4546
$ws.Range("K99").Value = 1139.2858
$ws.Range("L99").Value = 3202.4546
$ws.Range("M99").Value = 358.7141999999999
$ws.Range("N99").Value = -6198.4546

$ws.Range("H134").Value = 2160.875
$ws.Range("I134").Value = 2160.875
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6482.625
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3947.625
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2134.9
$ws.Range("I31").Value = 2677.875
$ws.Range("J31").Value = 1772.9166
$ws.Range("K31").Value = 2677.875
$ws.Range("L31").Value = 1772.9166
$ws.Range("M31").Value = -2382.875
$ws.Range("N31").Value = -2362.9166

$ws.Range("H34").Value = 2134.9
$ws.Range("I34").Value = 2677.875
$ws.Range("J34").Value = 1772.9166
$ws.Range("K34").Value = 2677.875
$ws.Range("L34").Value = 1772.9166
$ws.Range("M34").Value = -2475.875
$ws.Range("N34").Value = -2176.9166

$ws.Range("H134").Value = 2585
$ws.Range("I134").Value = 986.913
$ws.Range("J134").Value = 5926.4546
$ws.Range("K134").Value = 2960.739
$ws.Range("L134").Value = 17779.3638
$ws.Range("M134").Value = -425.739
$ws.Range("N134").Value = -22849.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1024.7046
$ws.Range("I131").Value = 742
$ws.Range("J131").Value = 1060.9487
$ws.Range("K131").Value = 2226
$ws.Range("L131").Value = 3182.8461
$ws.Range("M131").Value = 2814
$ws.Range("N131").Value = -13262.8461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9097.706
$ws.Range("I113").Value = 1978.2727
$ws.Range("J113").Value = 22150
$ws.Range("K113").Value = 1978.2727
$ws.Range("L113").Value = 22150
$ws.Range("M113").Value = 191.7273
$ws.Range("N113").Value = -26490

$ws.Range("H132").Value = 4787.4614
$ws.Range("I132").Value = 4835.3335
$ws.Range("J132").Value = 4213
$ws.Range("K132").Value = 14506.0005
$ws.Range("L132").Value = 12639
$ws.Range("M132").Value = -11976.0005
$ws.Range("N132").Value = -17699

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15086.571
$ws.Range("I61").Value = 15086.571
$ws.Range("K61").Value = 15086.571
$ws.Range("M61").Value = -14884.571

$ws.Range("H113").Value = 15086.571
$ws.Range("I113").Value = 15086.571
$ws.Range("K113").Value = 15086.571
$ws.Range("M113").Value = -12916.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 275.41666
$ws.Range("I113").Value = 227.81818
$ws.Range("K113").Value = 683.4545400000001
$ws.Range("M113").Value = 1486.54546

$ws.Range("H132").Value = 2698.5625
$ws.Range("I132").Value = 1909.25
$ws.Range("J132").Value = 3487.875
$ws.Range("K132").Value = 5727.75
$ws.Range("L132").Value = 10463.625
$ws.Range("M132").Value = -3197.75
$ws.Range("N132").Value = -15523.625
